$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new "MFR" / Replaceable-status column (H) to the BOM sheet.
# ---------------------------------------------------------------------------

# Step 1: seed column H (rows 2-61) from column G so every row inherits the
# same per-row cell style (border/alignment/number-format) already used by
# the sheet, instead of Excel's generic default style.
$ws.Range("G2:G61").Copy($ws.Range("H2:H61"))

# Step 2: row 3 is a "top of a bordered block" row in column G (style with a
# slightly different top border) but in column H it must use the regular
# style used throughout the rest of the column - copy that style from H4 and
# clear out the value.
$ws.Range("G4").Copy($ws.Range("H3"))
$ws.Range("H3").Value2 = ""

# Step 3: fill in the "potentially replaceable" rows first so that string
# lands first in the workbook's shared-string table.
foreach ($addr in @("H4","H5","H6","H7","H8","H9","H10","H23","H24","H25","H26","H27","H28","H29","H30","H31","H32","H33","H34","H35","H36","H37","H38","H39","H40","H41","H42","H43","H46","H47","H48","H49","H50","H53","H55","H56","H57","H58","H60","H61")) {
  $ws.Range($addr).Value2 = "potentially replaceable"
}

# Step 4: fill in the "not potentially replaceable" rows second.
foreach ($addr in @("H11","H12","H13","H14","H15","H16","H17","H18","H19","H20","H21","H22","H51","H52","H54")) {
  $ws.Range($addr).Value2 = "not potentially replaceable"
}

# Step 5: header label for the new column, added last.
$ws.Range("H2").Value2 = "Replaceable"

# Step 6: three existing MFR cells (G13, G38, G51) get re-styled with a
# dedicated Text-formatted (numFmtId 49) Courier New font (a distinct font
# entry carrying an explicit font-family classification) so the MFR codes
# stored there (some of which look numeric) are kept/shown as plain text.
foreach ($addr in @("G13","G38","G51")) {
  $rng = $ws.Range($addr)
  $rng.NumberFormat = "@"
  $rng.Font.Name = "Courier New"
  $rng.Font.Size = 10
  $rng.Font.Family = 3
}

# Step 7: size the new column to fit its (now longest) content, "not
# potentially replaceable".
$ws.Columns.Item(8).AutoFit()
$ws.Columns.Item(8).ColumnWidth = 27.15

# Step 8: reproduce the selection/scroll state left behind after the edit -
# the whole of column G selected (as if the user right-clicked the column
# header while adding the new column next to it).
$ws.Columns.Item(7).Select()
